$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 2098.5454
$ws.Range("I96").Value = 2120
$ws.Range("J96").Value = 2086.2856
$ws.Range("K96").Value = 6360
$ws.Range("L96").Value = 6258.8568
$ws.Range("M96").Value = -4987
$ws.Range("N96").Value = -9004.856800000001
$ws.Range("H137").Value = 6500.222
$ws.Range("I137").Value = 8283.666999999999
$ws.Range("J137").Value = 2933.3333
$ws.Range("K137").Value = 24851.001
$ws.Range("L137").Value = 8799.999899999999
$ws.Range("M137").Value = -22301.001
$ws.Range("N137").Value = -13899.9999
$ws.Range("H138").Value = 193447.89
$ws.Range("I138").Value = 30197
$ws.Range("K138").Value = 90591
$ws.Range("M138").Value = -85451

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 66700
$ws.Range("J24").Value = 66700
$ws.Range("L24").Value = 66700
$ws.Range("N24").Value = -67448
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H68").Value = 100000
$ws.Range("J68").Value = 100000
$ws.Range("L68").Value = 100000
$ws.Range("N68").Value = -101622
$ws.Range("H71").Value = 100000
$ws.Range("J71").Value = 100000
$ws.Range("L71").Value = 300000
$ws.Range("N71").Value = -308112
$ws.Range("H80").Value = 23709.1
$ws.Range("J80").Value = 21110.111
$ws.Range("L80").Value = 21110.111
$ws.Range("N80").Value = -23106.111
$ws.Range("H81").Value = 100000
$ws.Range("J81").Value = 100000
$ws.Range("L81").Value = 100000
$ws.Range("N81").Value = -101996
$ws.Range("H82").Value = 24090.5
$ws.Range("J82").Value = 24090.5
$ws.Range("L82").Value = 24090.5
$ws.Range("N82").Value = -24812.5
$ws.Range("H83").Value = 23709.1
$ws.Range("J83").Value = 21110.111
$ws.Range("L83").Value = 63330.333
$ws.Range("N83").Value = -73314.333
$ws.Range("H84").Value = 100000
$ws.Range("J84").Value = 100000
$ws.Range("L84").Value = 300000
$ws.Range("N84").Value = -309984
$ws.Range("H85").Value = 24090.5
$ws.Range("J85").Value = 24090.5
$ws.Range("L85").Value = 24090.5
$ws.Range("N85").Value = -26586.5
$ws.Range("H86").Value = 33362266
$ws.Range("J86").Value = 33362266
$ws.Range("L86").Value = 33362266
$ws.Range("N86").Value = -33364638
$ws.Range("H89").Value = 33362266
$ws.Range("J89").Value = 33362266
$ws.Range("L89").Value = 100086798
$ws.Range("N89").Value = -100098654
$ws.Range("H100").Value = 66700
$ws.Range("J100").Value = 66700
$ws.Range("L100").Value = 66700
$ws.Range("N100").Value = -68864

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5348.294
$ws.Range("I31").Value = 1212.2222
$ws.Range("K31").Value = 1212.2222
$ws.Range("M31").Value = -917.2221999999999
$ws.Range("H34").Value = 5348.294
$ws.Range("I34").Value = 1212.2222
$ws.Range("K34").Value = 1212.2222
$ws.Range("M34").Value = -1010.2222

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1033.9131
$ws.Range("J5").Value = 2695
$ws.Range("L5").Value = 8085
$ws.Range("N5").Value = -8309
$ws.Range("H31").Value = 6000
$ws.Range("J31").Value = 6000
$ws.Range("L31").Value = 18000
$ws.Range("N31").Value = -18576
$ws.Range("H39").Value = 1533.0769
$ws.Range("J39").Value = 1888
$ws.Range("L39").Value = 5664
$ws.Range("N39").Value = -6252
$ws.Range("H60").Value = 2711
$ws.Range("I60").Value = 467.5
$ws.Range("J60").Value = 2924.6667
$ws.Range("K60").Value = 1402.5
$ws.Range("L60").Value = 8774.000100000001
$ws.Range("M60").Value = -1151.5
$ws.Range("N60").Value = -9276.000100000001
$ws.Range("H131").Value = 954.1842
$ws.Range("I131").Value = 307.5
$ws.Range("J131").Value = 1030.2646
$ws.Range("K131").Value = 922.5
$ws.Range("L131").Value = 3090.7938
$ws.Range("M131").Value = 4117.5
$ws.Range("N131").Value = -13170.7938
$ws.Range("H132").Value = 3714.2334
$ws.Range("I132").Value = 2875.238
$ws.Range("K132").Value = 25877.142
$ws.Range("M132").Value = -23347.142
$ws.Range("H134").Value = 5005.448
$ws.Range("I134").Value = 2627.0588
$ws.Range("J134").Value = 8374.833000000001
$ws.Range("K134").Value = 7881.176399999999
$ws.Range("L134").Value = 25124.499
$ws.Range("M134").Value = -2811.176399999999
$ws.Range("N134").Value = -35264.499
$ws.Range("H135").Value = 1033.9131
$ws.Range("J135").Value = 2695
$ws.Range("L135").Value = 24255
$ws.Range("N135").Value = -29325
$ws.Range("H136").Value = 1664.9166
$ws.Range("I136").Value = 1197.9
$ws.Range("K136").Value = 3593.7
$ws.Range("M136").Value = 1506.3
$ws.Range("H137").Value = 9812
$ws.Range("I137").Value = 12539.9
$ws.Range("J137").Value = 2992.25
$ws.Range("K137").Value = 37619.7
$ws.Range("L137").Value = 8976.75
$ws.Range("M137").Value = -32519.7
$ws.Range("N137").Value = -19176.75
$ws.Range("H139").Value = 3808.6875
$ws.Range("I139").Value = 1941.4615
$ws.Range("K139").Value = 5824.3845
$ws.Range("M139").Value = -684.3845000000001
$ws.Range("H140").Value = 2012.0834
$ws.Range("I140").Value = 1395.625
$ws.Range("J140").Value = 3245
$ws.Range("K140").Value = 4186.875
$ws.Range("L140").Value = 9735
$ws.Range("M140").Value = 993.125
$ws.Range("N140").Value = -20095
$ws.Range("H141").Value = 7117.25
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3224.2856
$ws.Range("I100").Value = 3116
$ws.Range("K100").Value = 3116
$ws.Range("M100").Value = -2575

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 50000
$ws.Range("J104").Value = 50000
$ws.Range("L104").Value = 50000
$ws.Range("N104").Value = -56988
$ws.Range("H123").Value = 27968.125
$ws.Range("J123").Value = 27968.125
$ws.Range("L123").Value = 27968.125
$ws.Range("N123").Value = -37768.125

Write-Host "Edit complete"